# "Generate Report for Archive"
#
# The localization status report is regenerated: the zh-cn / de-de status
# moves on from "Ready for handoff" to "In Translation", and the columns
# that display that status are re-sized to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: columns E (zh-cn) and F (de-de) hold the per-locale status.
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# Locale sheets: column C holds the "Status" value for that locale's table.
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# Re-fit the status columns now that the text is shorter than before
# ("Ready for handoff" -> "In Translation"), shrinking them from
# ~17.22 characters wide down to ~13.41 characters wide.
$newWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
